$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list: Price (D) and Volume(1h) (E) columns
# D-column values that look like plain decimals need NumberFormat "@" first
# so Excel COM stores them as text (matching the original inline-string cells)
# instead of silently converting them to numbers.

$ws.Range("D2").Value = "47.180.48"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "2.481.82"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.07"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.41"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.35"
$ws.Range("E10").Value = "  +3.56%  "
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.08"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").Value = "2.871.48"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "2.494.29"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "47.100.26"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.74"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.59"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("E22").Value = "  +13.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.15"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "244.69"
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.55"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.61"
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.27"
$ws.Range("E31").Value = "  -2.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.34"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.95"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.61"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.55"
$ws.Range("E40").Value = "  +4.90%  "
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.35"
$ws.Range("E43").Value = "  -4.17%  "
$ws.Range("E44").Value = "  -1.24%  "
$ws.Range("D45").Value = "1.980.21"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("E47").Value = "  -6.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.01"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("E49").Value = "  -2.42%  "
$ws.Range("E50").Value = "  -5.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.38"
$ws.Range("E51").Value = "  +2.63%  "
